# edit.ps1 - applies the "Add files via upload" revision to Installation Guide.docx
#
# Summary of changes:
#   1. "Geom Alg Palette 2021Jan"  -> "Geom Alg Palette 2021Mar"   (1st list item, quoted name)
#   2. "right hand" -> "right-hand" (drops the gramStart/gramEnd proofErr pair)
#   3. "Place you cursor" -> "Place your cursor"
#   4. "GeomAlg2021Jan src" (2nd list item)   -> "GeomAlg2021Mar src"
#   5. "GeomAlg2021Jan src" (Files: section)  -> "GeomAlg2021Mar src"
#   6. "Geom Alg Palette src" (Files: section, bold) -> "GeomAlg2021Mar Palette"
#   7. Remove the whole "Improvements List - ..." paragraph plus the blank
#      paragraph that immediately followed it.

$d = $word.ActiveDocument

# --- 1. "...Palette 2021Jan"." -> "...Palette 2021Mar"." --------------------
$d.Content.Find.Execute("Palette 2021Jan", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Palette 2021Mar", 2) | Out-Null

# --- 2. "right hand" -> "right-hand" ----------------------------------------
$d.Content.Find.Execute("right hand", $true, $false, $false, $false, $false, `
    $true, 1, $false, "right-hand", 2) | Out-Null

# --- 3. "Place you cursor" -> "Place your cursor" ---------------------------
$d.Content.Find.Execute("Place you cursor", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Place your cursor", 2) | Out-Null

# --- 4. "notebook GeomAlg2021Jan src." -> "notebook GeomAlg2021Mar src." ----
$d.Content.Find.Execute("notebook GeomAlg2021Jan src.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "notebook GeomAlg2021Mar src.", 2) | Out-Null

# --- 5. "GeomAlg2021Jan src " (bold, Files: list) -> "GeomAlg2021Mar src " --
$d.Content.Find.Execute("GeomAlg2021Jan src", $true, $false, $false, $false, $false, `
    $true, 1, $false, "GeomAlg2021Mar src", 2) | Out-Null

# --- 6. "Geom Alg Palette src" (bold, Files: list) -> "GeomAlg2021Mar Palette"
$d.Content.Find.Execute("Geom Alg Palette src", $true, $false, $false, $false, $false, `
    $true, 1, $false, "GeomAlg2021Mar Palette", 2) | Out-Null

# --- 7. Remove the "Improvements List" paragraph + the blank line after it -
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Improvements List*") {
        $next = $d.Paragraphs.Item($i + 1)
        $deleteRange = $d.Range($para.Range.Start, $next.Range.End)
        $deleteRange.Delete()
        break
    }
}
